$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.422.75'
$ws.Range("E2").Value = '  -0.07%  '
$ws.Range("D3").Value = '3.422.08'
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''584.75'
$ws.Range("E5").Value = '  +0.12%  '
$ws.Range("D6").Value = '''179.22'
$ws.Range("E6").Value = '  +1.53%  '
$ws.Range("D7").Value = '''0.623'
$ws.Range("E7").Value = '  +3.69%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '3.420.58'
$ws.Range("E9").Value = '  -0.53%  '
$ws.Range("E10").Value = '  +0.27%  '
$ws.Range("E11").Value = '  +1.13%  '
$ws.Range("D12").Value = '''0.414'
$ws.Range("E12").Value = '  -0.66%  '
$ws.Range("D13").Value = '4.018.58'
$ws.Range("E13").Value = '  -0.62%  '
$ws.Range("E14").Value = '  +0.75%  '
$ws.Range("D15").Value = '''29.51'
$ws.Range("E15").Value = '  -2.45%  '
$ws.Range("D16").Value = '66.453.65'
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("E17").Value = '  +0.48%  '
$ws.Range("D18").Value = '3.426.35'
$ws.Range("E18").Value = '  -0.91%  '
$ws.Range("E19").Value = '  -0.65%  '
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("D21").Value = '''367.70'
$ws.Range("E21").Value = '  -3.03%  '
$ws.Range("E22").Value = '  -2.61%  '
$ws.Range("E23").Value = '  +0.81%  '
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").Value = '''0.998'
$ws.Range("E24").Value = '  -0.18%  '
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").Value = '''0.0000127'
$ws.Range("E25").Value = '  +5.71%  '
$ws.Range("E26").Value = '  -0.39%  '
$ws.Range("D27").Value = '''9.84'
$ws.Range("E27").Value = '  +0.73%  '
$ws.Range("E28").Value = '  +1.93%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("D30").Value = '''5.80'
$ws.Range("E30").Value = '  -0.88%  '
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("D32").Value = '''23.46'
$ws.Range("E32").Value = '  -3.53%  '
$ws.Range("D33").Value = '''0.999'
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").Value = '''7.04'
$ws.Range("E34").Value = '  -1.99%  '
$ws.Range("E35").Value = '  -3.90%  '
$ws.Range("E36").Value = '  -1.72%  '
$ws.Range("D37").Value = '''162.80'
$ws.Range("E37").Value = '  +1.25%  '
$ws.Range("D38").Value = '''0.876'
$ws.Range("E38").Value = '  -1.51%  '
$ws.Range("D39").Value = '''27.72'
$ws.Range("E39").Value = '  -5.64%  '
$ws.Range("E40").Value = '  +0.38%  '
$ws.Range("D41").Value = '''2.59'
$ws.Range("E41").Value = '  -1.49%  '
$ws.Range("E42").Value = '  -0.93%  '
$ws.Range("D43").Value = '2.707.76'
$ws.Range("E43").Value = '  -0.88%  '
$ws.Range("D44").Value = '''6.32'
$ws.Range("E44").Value = '  -0.91%  '
$ws.Range("E45").Value = '  -0.80%  '
$ws.Range("E46").Value = '  +3.06%  '
$ws.Range("D47").Value = '''40.01'
$ws.Range("E47").Value = '  -1.32%  '
$ws.Range("D48").Value = '''334.15'
$ws.Range("E48").Value = '  +8.10%  '
$ws.Range("E49").Value = '  -2.57%  '
$ws.Range("E50").Value = '  +2.41%  '
$ws.Range("D51").Value = '''32.04'
$ws.Range("E51").Value = '  +4.44%  '
